# SuppXLS/Scen_ELC_RES99.xlsx update:
# Add a new "P*OIL*" fuel-filter scenario row to the RNW/RES-penetration
# table, inserted right after the P*PEA* row (between the existing
# P*PEA* and P*DIS* rows), pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12 - this shifts old rows 12-19
# (P*DIS*, P*HFO*, P*HYD*, P*WIN*, P*SOL*, P*BIO*, P*GEO*, P*OCE*)
# down to rows 13-20, and copies formatting/formulas down with them.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new "P*OIL*" scenario,
# following the same pattern as its sibling rows (9-14): no RNW flag
# in column A, fuel filter in F, and the standard UC columns J:M plus
# the per-period N:T formulas referencing the RNW-Level row (row 2).
$ws.Range("F12").Value = "P*OIL*"
$ws.Range("J12").Value = "ELCC,ELCD"
$ws.Range("K12").Value = "UC_FLO"
$ws.Range("L12").Value = "O"
$ws.Range("M12").Value = "UP"
$ws.Range("N12").Formula = "=IF(`$A12=1,C`$2-1,C`$2)"
$ws.Range("O12").Formula = "=IF(`$A12=1,D`$2-1,D`$2)"
$ws.Range("P12").Formula = "=IF(`$A12=1,E`$2-1,E`$2)"
$ws.Range("Q12").Formula = "=IF(`$A12=1,F`$2-1,F`$2)"
$ws.Range("R12").Formula = "=IF(`$A12=1,G`$2-1,G`$2)"
$ws.Range("S12").Formula = "=IF(`$A12=1,H`$2-1,H`$2)"
$ws.Range("T12").Formula = "=IF(`$A12=1,I`$2-1,I`$2)"

# Leave the selection where the author's session ended up.
$ws.Range("J29:J30").Select()
